# Budsjett.xlsx edit: "Lagt til switcher i budsjettet" (Added switches to the budget)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new rows at 15-16 for the new "Switch" / "Multilayer
#    switch" line items. Excel shifts everything below down by two
#    rows and copies formatting from the row immediately above the
#    insertion point (row 14), which already carries the styles we
#    need for the new line-item rows.
# ---------------------------------------------------------------------
$ws.Rows("15:16").Insert()

# ---------------------------------------------------------------------
# 2. Donasjoner til open source (row 5): quantity/unit price increase.
#    Copy the number format used for the other top "amount" cells
#    (B3/B4, style used for F5 in the target workbook) onto F5, then
#    set the new values.
# ---------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("F5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 10000
$ws.Range("G5").Formula = "=SUM(E5*F5)"

# Delsum; programvare subtotal (row 6) becomes a real formula instead
# of a hard-coded value.
$ws.Range("G6").Formula = "=SUM(G5)"

# ---------------------------------------------------------------------
# 3. New row 15 : "Switch"
# ---------------------------------------------------------------------
$ws.Range("D15").Value = "Switch"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 951
$ws.Range("G15").Formula = "=SUM(E15*F15)"

# New row 16 : "Multilayer switch"
$ws.Range("D16").Value = "Multilayer switch"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 2544
$ws.Range("G16").Formula = "=SUM(E16*F16)"

# Hyperlinks for the two new hardware line items (same pattern as the
# existing hardware rows D10/D13/D14).
$ws.Hyperlinks.Add($ws.Range("D15"), "https://www.komplett.no/product/1251402/nettverk-lagring/nettverk/switcher/netgear-gs308-unmanaged-switch") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D16"), "https://www.komplett.no/product/1253751/nettverk-lagring/nettverk/switcher/netgear-gs724tv4-managed-switch") | Out-Null

# ---------------------------------------------------------------------
# 4. Delsum; maskinvare subtotal (now row 17) must sum through the new
#    rows as well.
# ---------------------------------------------------------------------
$ws.Range("G17").Formula = "=SUM(G10:G16)"

Write-Output "Budget updated: switches added."
